$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.729.85'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '1.974.31'
$ws.Range("E3").Value = '  -2.83%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.14'
$ws.Range("E5").Value = '  -7.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.598'
$ws.Range("E6").Value = '  -3.87%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.49'
$ws.Range("E8").Value = '  -7.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.367'
$ws.Range("E9").Value = '  -4.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.75'
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0740'
$ws.Range("E11").Value = '  -5.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0975'
$ws.Range("E12").Value = '  -4.03%  '
$ws.Range("D13").Value = '2.273.41'
$ws.Range("E13").Value = '  -2.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '13.83'
$ws.Range("E14").Value = '  -5.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.21'
$ws.Range("E15").Value = '  -5.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.744'
$ws.Range("E16").Value = '  -9.51%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.988.52'
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.98'
$ws.Range("E18").Value = '  -7.26%  '
$ws.Range("D19").Value = '36.725.19'
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.43'
$ws.Range("E20").Value = '  -3.19%  '
$ws.Range("D21").Value = '0.0₃0795'
$ws.Range("E21").Value = '  -6.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.00'
$ws.Range("E22").Value = '  -4.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '220.34'
$ws.Range("E23").Value = '  -3.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.35'
$ws.Range("E26").Value = '  -9.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.06'
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.50'
$ws.Range("E28").Value = '  -6.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.90'
$ws.Range("E29").Value = '  -5.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.123'
$ws.Range("E30").Value = '  -5.54%  '
$ws.Range("E31").Value = '  -5.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.116'
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.34'
$ws.Range("E33").Value = '  -8.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0601'
$ws.Range("E34").Value = '  -9.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.17'
$ws.Range("E35").Value = '  -9.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.27'
$ws.Range("E36").Value = '  -6.99%  '
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("E38").Value = '  -1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.20'
$ws.Range("E39").Value = '  -5.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.19'
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("D42").Value = '1.413.54'
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0201'
$ws.Range("E43").Value = '  -6.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0877'
$ws.Range("E44").Value = '  -9.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.10'
$ws.Range("E45").Value = '  -7.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '86.71'
$ws.Range("E46").Value = '  -4.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.988'
$ws.Range("E47").Value = '  -5.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '14.82'
$ws.Range("E48").Value = '  -8.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.87'
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.65'
$ws.Range("E50").Value = '  -9.88%  '
$ws.Range("D51").Value = '2.165.41'
$ws.Range("E51").Value = '  -2.38%  '

Write-Host "Update complete"